$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.519.40'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '1.854.37'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6333'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07556'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2982'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07737'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '1.882.64'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.023'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6920'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009911'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '2.134.30'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.289'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.91%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '29.560.04'
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '233.90'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.690'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.68%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9996'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1400'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.482'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.94%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.480'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05942'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.37%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.251'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.31%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.132'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.36%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.034'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.907'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.172'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7245'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.79%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.584'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").Value = '1.242.33'
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.800'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01797'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9084'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.111'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").Value = '2.044.72'
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9991'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.401'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4051'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.142'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000117'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.12%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.712'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.24%  '
